$d = $word.ActiveDocument

$found = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "https://www.youtube.com/watch?v=XtR1P4BfuQQ") {
        $p.Range.Font.Color = 5880731
        $found = $true
        break
    }
}

if (-not $found) {
    throw "Target paragraph not found"
}
